$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1761.7333
$ws.Range("I40").Value = 2205.75
$ws.Range("J40").Value = 1465.7222
$ws.Range("K40").Value = 2205.75
$ws.Range("L40").Value = 1465.7222
$ws.Range("M40").Value = -2030.75
$ws.Range("N40").Value = -1815.7222
$ws.Range("H116").Value = 4434.3335
$ws.Range("I116").Value = 1200
$ws.Range("J116").Value = 5081.2
$ws.Range("K116").Value = 1200
$ws.Range("L116").Value = 5081.2
$ws.Range("M116").Value = 2242
$ws.Range("N116").Value = -11965.2
$ws.Range("H132").Value = 2269487.8
$ws.Range("I132").Value = 2599295
$ws.Range("J132").Value = 2062.75
$ws.Range("K132").Value = 7797885
$ws.Range("L132").Value = 6188.25
$ws.Range("M132").Value = -7795355
$ws.Range("N132").Value = -11248.25
$ws.Range("H137").Value = 1956.3889
$ws.Range("I137").Value = 1730.7142
$ws.Range("J137").Value = 2100
$ws.Range("K137").Value = 5192.142599999999
$ws.Range("L137").Value = 6300
$ws.Range("M137").Value = -2642.142599999999
$ws.Range("N137").Value = -11400

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3715
$ws.Range("I32").Value = 3256.1042
$ws.Range("J32").Value = 14728.5
$ws.Range("K32").Value = 3256.1042
$ws.Range("L32").Value = 14728.5
$ws.Range("M32").Value = -2969.1042
$ws.Range("N32").Value = -15302.5
$ws.Range("H45").Value = 30303892
$ws.Range("I45").Value = 47619784
$ws.Range("J45").Value = 1083.5
$ws.Range("K45").Value = 47619784
$ws.Range("L45").Value = 1083.5
$ws.Range("M45").Value = -47619407
$ws.Range("N45").Value = -1837.5
$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678
$ws.Range("H132").Value = 4283.4556
$ws.Range("I132").Value = 4764.9653
$ws.Range("J132").Value = 2953.5715
$ws.Range("K132").Value = 14294.8959
$ws.Range("L132").Value = 8860.7145
$ws.Range("M132").Value = -11764.8959
$ws.Range("N132").Value = -13920.7145

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 449.30768
$ws.Range("J64").Value = 460.875
$ws.Range("L64").Value = 460.875
$ws.Range("N64").Value = -910.875
$ws.Range("H67").Value = 449.30768
$ws.Range("J67").Value = 460.875
$ws.Range("L67").Value = 460.875
$ws.Range("N67").Value = -2020.875
$ws.Range("H86").Value = 2678.889
$ws.Range("I86").Value = 2738.25
$ws.Range("J86").Value = 2631.4
$ws.Range("K86").Value = 2738.25
$ws.Range("L86").Value = 2631.4
$ws.Range("M86").Value = -1615.25
$ws.Range("N86").Value = -4877.4
$ws.Range("H89").Value = 2678.889
$ws.Range("I89").Value = 2738.25
$ws.Range("J89").Value = 2631.4
$ws.Range("K89").Value = 13691.25
$ws.Range("L89").Value = 13157
$ws.Range("M89").Value = -8075.25
$ws.Range("N89").Value = -24389
$ws.Range("H134").Value = 6625.6924
$ws.Range("I134").Value = 8456.235000000001
$ws.Range("J134").Value = 3168
$ws.Range("K134").Value = 25368.705
$ws.Range("L134").Value = 9504
$ws.Range("M134").Value = -22833.705
$ws.Range("N134").Value = -14574

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3690.818
$ws.Range("I62").Value = 2735
$ws.Range("J62").Value = 4837.8
$ws.Range("K62").Value = 2735
$ws.Range("L62").Value = 4837.8
$ws.Range("M62").Value = -2111
$ws.Range("N62").Value = -6085.8
$ws.Range("H65").Value = 3690.818
$ws.Range("I65").Value = 2735
$ws.Range("J65").Value = 4837.8
$ws.Range("K65").Value = 13675
$ws.Range("L65").Value = 24189
$ws.Range("M65").Value = -10555
$ws.Range("N65").Value = -30429
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20663.2
$ws.Range("J131").Value = 27844.783
$ws.Range("L131").Value = 83534.349
$ws.Range("N131").Value = -93614.349

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6187.857
$ws.Range("I80").Value = 4219.1665
$ws.Range("J80").Value = 18000
$ws.Range("K80").Value = 4219.1665
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -3221.1665
$ws.Range("N80").Value = -19996
$ws.Range("H83").Value = 6187.857
$ws.Range("I83").Value = 4219.1665
$ws.Range("J83").Value = 18000
$ws.Range("K83").Value = 21095.8325
$ws.Range("L83").Value = 90000
$ws.Range("M83").Value = -16103.8325
$ws.Range("N83").Value = -99984
$ws.Range("H102").Value = 1265.7097
$ws.Range("I102").Value = 970.6539
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 970.6539
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = 651.3461
$ws.Range("N102").Value = -6044
$ws.Range("H123").Value = 28072.363
$ws.Range("J123").Value = 28072.363
$ws.Range("L123").Value = 28072.363
$ws.Range("N123").Value = -32972.363
$ws.Range("H132").Value = 3640.4666
$ws.Range("I132").Value = 3755.3
$ws.Range("J132").Value = 2721.8
$ws.Range("K132").Value = 11265.9
$ws.Range("L132").Value = 8165.400000000001
$ws.Range("M132").Value = -8735.900000000001
$ws.Range("N132").Value = -13225.4
$ws.Range("H134").Value = 12949.4
$ws.Range("J134").Value = 12949.4
$ws.Range("L134").Value = 38848.2
$ws.Range("N134").Value = -43918.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 18333.334
$ws.Range("J110").Value = 18333.334
$ws.Range("L110").Value = 18333.334
$ws.Range("N110").Value = -26513.334
$ws.Range("H122").Value = 5746.1816
$ws.Range("I122").Value = 9421.6
$ws.Range("K122").Value = 28264.8
$ws.Range("M122").Value = -25814.8
$ws.Range("H132").Value = 11385.772
$ws.Range("I132").Value = 21439.7
$ws.Range("J132").Value = 3007.5
$ws.Range("K132").Value = 64319.10000000001
$ws.Range("L132").Value = 9022.5
$ws.Range("M132").Value = -61789.10000000001
$ws.Range("N132").Value = -14082.5
$ws.Range("H136").Value = 10236.286
$ws.Range("I136").Value = 14067.556
$ws.Range("J136").Value = 3340
$ws.Range("K136").Value = 42202.66800000001
$ws.Range("L136").Value = 10020
$ws.Range("M136").Value = -39652.66800000001
$ws.Range("N136").Value = -15120

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 34860
$ws.Range("J119").Value = 47290
$ws.Range("L119").Value = 47290
$ws.Range("N119").Value = -56966
$ws.Range("H132").Value = 1652.5714
$ws.Range("I132").Value = 1451.6487
$ws.Range("J132").Value = 3139.4
$ws.Range("K132").Value = 4354.9461
$ws.Range("L132").Value = 9418.200000000001
$ws.Range("M132").Value = -1824.9461
$ws.Range("N132").Value = -14478.2
